$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.498.27"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").Value = "1.671.98"
$ws.Range("E3").Value = "  +1.62%  "

$ws.Range("D5").Value = "220.54"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").Value = "0.5278"
$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  +2.72%  "

$ws.Range("D9").Value = "0.06378"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "21.77"
$ws.Range("E10").Value = "  +4.87%  "

$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.673.80"

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.488"
$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").Value = "0.5569"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").Value = "0.0₅8313"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("D16").Value = "65.59"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").Value = "26.501.09"
$ws.Range("E17").Value = "  +1.89%  "

$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").Value = "4.760"
$ws.Range("E19").Value = "  +1.10%  "

$ws.Range("D20").Value = "193.07"
$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("D21").Value = "10.33"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").Value = "6.305"
$ws.Range("E22").Value = "  +0.85%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "0.1269"
$ws.Range("E24").Value = "  +3.99%  "

$ws.Range("D25").Value = "139.43"
$ws.Range("E25").Value = "  -3.47%  "

$ws.Range("D26").Value = "7.404"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  +2.96%  "

$ws.Range("D28").Value = "1.425"
$ws.Range("E28").Value = "  +2.66%  "

$ws.Range("D29").Value = "0.06201"
$ws.Range("E29").Value = "  +4.42%  "

$ws.Range("D30").Value = "1.292"
$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("D31").Value = "3.611"
$ws.Range("E31").Value = "  +6.21%  "

$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("E34").Value = "  +1.04%  "

$ws.Range("D35").Value = "0.6081"
$ws.Range("E35").Value = "  +8.11%  "

$ws.Range("D36").Value = "2.415"
$ws.Range("E36").Value = "  +0.82%  "

$ws.Range("D37").Value = "2.778"
$ws.Range("E37").Value = "  +0.85%  "

$ws.Range("D38").Value = "0.01616"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("D39").Value = "6.039"
$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("D40").Value = "1.086.48"
$ws.Range("E40").Value = "  +5.57%  "

$ws.Range("D41").Value = "0.8577"
$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").Value = "100.58"
$ws.Range("E43").Value = "  +1.72%  "

$ws.Range("D44").Value = "1.816.78"
$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("D46").Value = "58.38"
$ws.Range("E46").Value = "  +4.89%  "

$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").Value = "1.522"
$ws.Range("E48").Value = "  +10.42%  "

$ws.Range("D49").Value = "8.110"
$ws.Range("E49").Value = "  +0.82%  "

$ws.Range("D50").Value = "0.05194"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").Value = "6.003"
$ws.Range("E51").Value = "  +1.71%  "
